# Update "想去人数" (F column) counts that changed between crawl runs.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 76
$ws1.Range("F4").Value = 2086
$ws1.Range("F5").Value = 368
$ws1.Range("F6").Value = 635
$ws1.Range("F9").Value = 10711
$ws1.Range("F12").Value = 287
$ws1.Range("F14").Value = 420
$ws1.Range("F15").Value = 7569
$ws1.Range("F16").Value = 1114
$ws1.Range("F18").Value = 263
$ws1.Range("F19").Value = 67
$ws1.Range("F20").Value = 3338

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 76
$ws4.Range("F4").Value = 2086
$ws4.Range("F5").Value = 368
$ws4.Range("F6").Value = 635
$ws4.Range("F7").Value = 24
$ws4.Range("F12").Value = 10711
$ws4.Range("F15").Value = 287
$ws4.Range("F17").Value = 420
$ws4.Range("F18").Value = 7569
$ws4.Range("F19").Value = 1114
$ws4.Range("F21").Value = 263
$ws4.Range("F22").Value = 67
$ws4.Range("F23").Value = 3338
